$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New collection rows (MCH336-1 .. MCH336-5)
$data = @(
    @("MCH336-1", "PUBLICATIONS- MOZAMBIQUE REVOLUTION, SOUTHERN AFRICA", "Series", "1 Box", "LOCATION: 33H | GRAP COUNT NUMER: NONE"),
    @("MCH336-2", "PRESS CLIPPINGS- NATIONAL PARTY, SADF/MURDER/POLITICAL VIOLENCE, UNIVERSITIES", "Series", "1 Box", "LOCATION: 33H | GRAP COUNT NUMER: NONE"),
    @("MCH336-3", "PUBLICATIONS- EAST AFRICA JOURNAL, MOZAMBIQUE REVOLUTION", "Series", "1 Box", "LOCATION: 33H | GRAP COUNT NUMER: NONE"),
    @("MCH336-4", "PUBLICATIONS- ANGOLA BULLETIN, BLATTER DES IZ3W, THIRD WORLD JOURNAL, INTELLIGENCE & FOREIGN POLICY, PRESERVATION & DEVELOPMENT OF INDIGENOUS ARTS, MOZAMBIQUE REVOLUTION (TO BE SORTED)", "Series", "1 Box", "LOCATION: 33H | GRAP COUNT NUMER: NONE"),
    @("MCH336-5", "PUBLICATIONS- BLACK ORPHEUS, EAST AFRICA JOURNAL ( TO BE SORTED)", "Series", "1 Box", "LOCATION: 33H | GRAP COUNT NUMER: NONE")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
}

# Apply the new-row font (Calibri 10, automatic/theme text color) to every
# touched column individually so column B (left untouched by the source
# data) is never materialised.
foreach ($col in @("A", "C", "D", "E", "F", "G", "H")) {
    $rng = $ws.Range("${col}2:${col}6")
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}
